# Update the "想去人数" (interested-count) values in column F across the
# four worksheets, per the latest site regeneration (gh-pages build
# 456a3b4). Each entry below maps a worksheet name + row number to its
# new column-F value.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        6  = 343
        7  = 1144
        9  = 7044
        11 = 88
        13 = 7928
        16 = 5481
        18 = 2383
        20 = 4551
        21 = 289
        22 = 382
        25 = 350
        26 = 249
        28 = 2242
        30 = 257
        32 = 108
        33 = 569
        36 = 1465
        37 = 30
        39 = 2250
        40 = 2202
        42 = 2
    }
    "演出" = @{
        5 = 21
    }
    "本地生活" = @{
        3 = 1273
    }
    "全部类型" = @{
        4  = 1273
        7  = 343
        8  = 1145
        10 = 7044
        12 = 88
        14 = 7928
        17 = 5481
        19 = 2383
        21 = 4551
        22 = 289
        23 = 382
        28 = 350
        29 = 249
        31 = 2242
        33 = 257
        35 = 108
        36 = 569
        39 = 21
        40 = 1465
        41 = 30
        43 = 2250
        45 = 2202
        47 = 2
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
